$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fields")

# Insert two new rows right after the header block for the "models" table
# (new rows 7 and 8), pushing all subsequent rows down by two.
$null = $ws.Rows.Item(7).Insert()
$null = $ws.Rows.Item(7).Insert()

# Table name column (reuses existing "models" shared string)
$ws.Cells.Item(7,1).Value = "models"
$ws.Cells.Item(8,1).Value = "models"

# Field names (order matters so shared strings are appended in the same
# sequence as the target workbook: US state field, then country field)
$ws.Cells.Item(7,2).Value = "geographic_resolution_US_state"
$ws.Cells.Item(8,2).Value = "geographic_resolution_country"

# Definitions (appended next, US state then country)
$ws.Cells.Item(7,3).Value = "Whether or not the model currently produced projections for one or more US states (e.g., are results available at the state level in the US?)"
$ws.Cells.Item(8,3).Value = "Whether or not the model currently produced projections for one or more countries (e.g., are results available at the country level, in the US or internationally?)"

# Possible values column (new shared string "Boolean", reused on both rows)
$ws.Cells.Item(7,4).Value = "Boolean"
$ws.Cells.Item(8,4).Value = "Boolean"

# Match row heights used in the published workbook
$ws.Rows.Item(7).RowHeight = 56
$ws.Rows.Item(8).RowHeight = 56

# Restore the sheet selection/scroll position
$null = $ws.Activate()
$null = $ws.Range("D9").Select()
